$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "5" values for the cells that changed from blank to 5
$ws.Range("F7").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E20").Value = 5
$ws.Range("F25").Value = 5
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 5

# Update the active selection to F25 (and scroll the frozen pane up to row 4)
$ws.Range("F25").Select()
